$d = $word.ActiveDocument

# 1) "fspe" -> "frane" (both occurrences)
$d.Content.Find.Execute("fspe", $false, $false, $false, $false, $false, $true, 1, $false, "frane", 2) | Out-Null

# 2) Split "Fletes SPERONI" into two runs: "Fletes " and "DEFENSA "
$d.Content.Find.Execute("Fletes SPERONI", $false, $false, $false, $false, $false, $true, 1, $false, "Fletes ", 2) | Out-Null

$fletesPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Fletes `r") {
        $fletesPara = $d.Paragraphs($i)
        break
    }
}
if ($fletesPara -eq $null) {
    throw "Could not locate the 'Fletes ' paragraph"
}
$fr = $fletesPara.Range
$splitPos = $fr.End - 1
$splitRng = $d.Range($splitPos, $splitPos)
$splitRng.InsertParagraphAfter()
$afterRng = $d.Range($splitPos + 1, $splitPos + 1)
$afterRng.InsertAfter("DEFENSA ")
$mergeRng = $d.Range($splitPos, $splitPos + 1)
$mergeRng.Delete()

# 3) Move the "_GoBack" bookmark onto the "Realizamos..." paragraph
$d.Bookmarks("_GoBack").Delete()

$mudanzaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Realizamos la mudanza*") {
        $mudanzaPara = $d.Paragraphs($i)
        break
    }
}
if ($mudanzaPara -eq $null) {
    throw "Could not locate the 'Realizamos la mudanza...' paragraph"
}
$mr = $mudanzaPara.Range
$bmPos = $mr.End - 1
$bmRng = $d.Range($bmPos, $bmPos)
$bmRng.InsertAfter("X")
$bmRng2 = $d.Range($bmPos, $bmPos + 1)
$bmRng2.Bookmarks.Add("_GoBack")
$bmRng2.Text = ""
